$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the old _GoBack bookmark (it will be re-added later at
#    its new location further down the document).
# ---------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------
# 2) Rewrite the "Dalam tahap matrikulasi ..." paragraph segment.
#    This single Find/Replace spans several original runs
#    (including the removed bookmark location) and replaces them
#    with the new wording.
# ---------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    " Dalam tahap matrikulasi ini mahasiswa diwajibkan mengikuti berbagai program (sebagai syarat lulus tahap matrikulasi) didalamnya yaitu program Pembinaan, program Akademik dan program Bahasa (TLC / Tazkia Language Center). Ketiga program tersebut haruslah di monitor dengan baik oleh pihak manajemen matrikuklasi agar nantinya data bisa diolah dengan baik hingga dapat dijadikan suatu informasi yang mudah dibaca oleh pihak berkepentingan. Pada kenyataannya, seluruh kegiatan pada program tersebut belum ada suatu sistem yang menangani",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Dalam tahap matrikulasi ini mahasiswa diwajibkan mengikuti berbagai kegiatan  didalamnya dan ditargetkan jumlah total presensi sesuai ketentuannya. Kegiatan tersebut diantaranya adalah shalat wajib berjamaah, tahsin/tahfidz dan ta’lim. Adapun target lain yang perlu dicapai oleh mahasiswa matrikulasi yaitu hafalan quran. Beberapa kegiatan tersebut haruslah di monitor dengan baik oleh pihak manajemen matrikuklasi agar nantinya data bisa diolah dengan baik hingga dapat dijadikan suatu informasi yang mudah dibaca oleh pihak berkepentingan. Pada kenyataannya, kegiatan-kegiatan tersebut belum ada suatu sistem yang menangani",
    2)
Write-Output $found
